# daily needs jsp added
# Adds a new "DailyNeedsListing" worksheet (after "PathLabsListing") populated
# with three daily-needs shop listings (Family Store, Mini Mart, One Stop Mart).

$wb = $excel.ActiveWorkbook

$pathLabs = $wb.Worksheets.Item("PathLabsListing")

# New sheet goes right after PathLabsListing, and becomes the active tab.
$ws = $wb.Worksheets.Add($null, $pathLabs)
$ws.Name = "DailyNeedsListing"

# Header row
$ws.Range("A1").Value = "number"
$ws.Range("B1").Value = "active"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "address"
$ws.Range("F1").Value = "contactNumber"
$ws.Range("G1").Value = "website"
$ws.Range("H1").Value = "openTime"
$ws.Range("I1").Value = "imageUrl"
$ws.Range("J1").Value = "map"
$ws.Range("K1").Value = "market"

# Row 2 - Family Store
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Y"
$ws.Range("C2").Value = "DailyNeeds"
$ws.Range("D2").Value = "Family Store"
$ws.Range("E2").Value = "LGF 41, City Plaza, Gaur City, Greater Noida West (Noida Extension)"
$ws.Range("F2").Value = "09599975791, 09599975792, 07827439297"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "07:30AM to 10:00PM"
$ws.Range("I2").Value = "/img/dailyNeeds/dailyNeeds_FamilyStore.jpeg"
$ws.Range("J2").Value = "https://www.google.com/maps/place/Family+Store/@28.6146087,77.4255353,17z/data=!3m1!4b1!4m5!3m4!1s0x390cee4b8783e553:0x652d6c4956ab52f6!8m2!3d28.614604!4d77.427724"
$ws.Range("K2").Value = "City Plaza, Gaur City"

# Row 3 - Mini Mart
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "DailyNeeds"
$ws.Range("D3").Value = "Mini Mart"
$ws.Range("E3").Value = "Shop no.23, City plaza, Gaur city 1, Greater Noida West (Noida Extension)"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "08800563608"
$ws.Range("G3").Value = "www.Minimart.com "
$ws.Range("H3").Value = "08:00AM to 09:30PM"
$ws.Range("I3").Value = "/img/dailyNeeds/dailyNeeds_MiniMart.jpeg"
$ws.Range("J3").Value = "https://www.google.com/maps/place/Mini+Mart/@28.6150307,77.4255613,17z/data=!3m1!4b1!4m5!3m4!1s0x390cee4b85c7e24d:0x984467dc816981a0!8m2!3d28.615026!4d77.42775"
$ws.Range("K3").Value = "City Plaza, Gaur City"

# Row 4 - One Stop Mart
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "DailyNeeds"
$ws.Range("D4").Value = "One Stop Mart"
$ws.Range("E4").Value = "Shop Number 45-50, Lower Ground Floor, Galliria Market, Gaur City 2, Greater Noida West (Noida Extension)"
$ws.Range("F4").Value = "07290016382, 07290016383"
$ws.Range("G4").Value = "N/A"
$ws.Range("H4").Value = "24 Hrs"
$ws.Range("I4").Value = "/img/dailyNeeds/dailyNeeds_OneStopMart.jpeg"
$ws.Range("J4").Value = "https://www.google.com/maps/place/One+Stop+Mart/@28.619435,77.4210851,17z/data=!3m1!4b1!4m5!3m4!1s0x390cefb559f85485:0x6e99725b0e217cc2!8m2!3d28.6194303!4d77.4232738"
$ws.Range("K4").Value = "Galliria Market, Gaur City"

# Column widths to match the other *Listing sheets
$ws.Columns.Item(1).ColumnWidth = 7.453125
$ws.Columns.Item(2).ColumnWidth = 5.7265625
$ws.Columns.Item(3).ColumnWidth = 10.08984375
$ws.Columns.Item(4).ColumnWidth = 21.453125
$ws.Columns.Item(5).ColumnWidth = 60
$ws.Columns.Item(6).ColumnWidth = 14
$ws.Columns.Item(7).ColumnWidth = 17.08984375
$ws.Columns.Item(8).ColumnWidth = 18.36328125
$ws.Columns.Item(9).ColumnWidth = 36.54296875
$ws.Columns.Item(10).ColumnWidth = 61.6328125
$ws.Columns.Item(11).ColumnWidth = 20.26953125

$ws.Range("E15").Select()
